# Update the date heading and all two-digit multiplication problems
# in the table to the new set of values.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-14 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-15 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("52×40=", $true, $false, $false, $false, $false, $true, 1, $false, "61×33=", 2) | Out-Null
$d.Content.Find.Execute("26×83=", $true, $false, $false, $false, $false, $true, 1, $false, "45×59=", 2) | Out-Null
$d.Content.Find.Execute("11×31=", $true, $false, $false, $false, $false, $true, 1, $false, "66×61=", 2) | Out-Null
$d.Content.Find.Execute("95×22=", $true, $false, $false, $false, $false, $true, 1, $false, "83×23=", 2) | Out-Null
$d.Content.Find.Execute("77×67=", $true, $false, $false, $false, $false, $true, 1, $false, "75×32=", 2) | Out-Null
$d.Content.Find.Execute("85×12=", $true, $false, $false, $false, $false, $true, 1, $false, "22×84=", 2) | Out-Null
$d.Content.Find.Execute("98×24=", $true, $false, $false, $false, $false, $true, 1, $false, "51×13=", 2) | Out-Null
$d.Content.Find.Execute("37×99=", $true, $false, $false, $false, $false, $true, 1, $false, "23×90=", 2) | Out-Null
$d.Content.Find.Execute("75×14=", $true, $false, $false, $false, $false, $true, 1, $false, "87×21=", 2) | Out-Null
$d.Content.Find.Execute("62×76=", $true, $false, $false, $false, $false, $true, 1, $false, "52×57=", 2) | Out-Null
$d.Content.Find.Execute("85×16=", $true, $false, $false, $false, $false, $true, 1, $false, "82×11=", 2) | Out-Null
$d.Content.Find.Execute("54×81=", $true, $false, $false, $false, $false, $true, 1, $false, "25×59=", 2) | Out-Null
$d.Content.Find.Execute("12×70=", $true, $false, $false, $false, $false, $true, 1, $false, "47×99=", 2) | Out-Null
$d.Content.Find.Execute("55×65=", $true, $false, $false, $false, $false, $true, 1, $false, "48×16=", 2) | Out-Null
$d.Content.Find.Execute("35×46=", $true, $false, $false, $false, $false, $true, 1, $false, "20×59=", 2) | Out-Null
$d.Content.Find.Execute("12×54=", $true, $false, $false, $false, $false, $true, 1, $false, "75×87=", 2) | Out-Null
$d.Content.Find.Execute("47×23=", $true, $false, $false, $false, $false, $true, 1, $false, "67×22=", 2) | Out-Null
$d.Content.Find.Execute("57×46=", $true, $false, $false, $false, $false, $true, 1, $false, "14×70=", 2) | Out-Null
$d.Content.Find.Execute("38×27=", $true, $false, $false, $false, $false, $true, 1, $false, "74×57=", 2) | Out-Null
$d.Content.Find.Execute("52×31=", $true, $false, $false, $false, $false, $true, 1, $false, "27×23=", 2) | Out-Null
$d.Content.Find.Execute("52×17=", $true, $false, $false, $false, $false, $true, 1, $false, "41×93=", 2) | Out-Null
$d.Content.Find.Execute("17×43=", $true, $false, $false, $false, $false, $true, 1, $false, "96×41=", 2) | Out-Null
$d.Content.Find.Execute("57×98=", $true, $false, $false, $false, $false, $true, 1, $false, "26×89=", 2) | Out-Null
$d.Content.Find.Execute("16×71=", $true, $false, $false, $false, $false, $true, 1, $false, "50×37=", 2) | Out-Null
$d.Content.Find.Execute("57×73=", $true, $false, $false, $false, $false, $true, 1, $false, "85×54=", 2) | Out-Null
